$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C5").Value = 25
